$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the 属性值 (stat value) in H9 from 708 to 793.
# The cell is formatted as text (numFmtId 49 => "@"), so assign it as a string
# to keep it a shared-string text cell; Excel will recalc the dependent
# formulas in H10:H13 automatically.
$ws.Range("H9").Value = "793"

# Move the active selection to H15 (last used cell), matching the saved
# sheetView state after editing.
$ws.Range("H15").Select()

$wb.Save()
